$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value2 = 3096.6667
$ws.Range("J46").Value2 = 5000
$ws.Range("L46").Value2 = 15000
$ws.Range("N46").Value2 = -15238
$ws.Range("H53").Value2 = 449.77777
$ws.Range("I53").Value2 = 234.38461
$ws.Range("J53").Value2 = 1009.8
$ws.Range("K53").Value2 = 234.38461
$ws.Range("L53").Value2 = 1009.8
$ws.Range("M53").Value2 = 402.61539
$ws.Range("N53").Value2 = -2283.8
$ws.Range("H60").Value2 = 3096.6667
$ws.Range("J60").Value2 = 5000
$ws.Range("L60").Value2 = 15000
$ws.Range("N60").Value2 = -15968
$ws.Range("H69").Value2 = 18733.2
$ws.Range("J69").Value2 = 19428.428
$ws.Range("L69").Value2 = 58285.284
$ws.Range("N69").Value2 = -60033.284
$ws.Range("H72").Value2 = 18733.2
$ws.Range("J72").Value2 = 19428.428
$ws.Range("L72").Value2 = 174855.852
$ws.Range("N72").Value2 = -183591.852
$ws.Range("H76").Value2 = 8347
$ws.Range("I76").Value2 = 0
$ws.Range("J76").Value2 = 8347
$ws.Range("K76").Value2 = 0
$ws.Range("L76").Value2 = 8347
$ws.Range("N76").Value2 = -8977
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value2 = 8347
$ws.Range("I79").Value2 = 0
$ws.Range("J79").Value2 = 8347
$ws.Range("K79").Value2 = 0
$ws.Range("L79").Value2 = 8347
$ws.Range("N79").Value2 = -10531
$ws.Range("M79").ClearContents()
$ws.Range("H128").Value2 = 0
$ws.Range("J128").Value2 = 0
$ws.Range("L128").Value2 = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value2 = 47500
$ws.Range("J130").Value2 = 47500
$ws.Range("L130").Value2 = 47500
$ws.Range("N130").Value2 = -57540
$ws.Range("H137").Value2 = 3498.5264
$ws.Range("I137").Value2 = 2826.4167
$ws.Range("K137").Value2 = 8479.250100000001
$ws.Range("M137").Value2 = -5929.250100000001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 3290.9167
$ws.Range("I74").Value2 = 1245.8572
$ws.Range("J74").Value2 = 17606.334
$ws.Range("K74").Value2 = 1245.8572
$ws.Range("L74").Value2 = 17606.334
$ws.Range("M74").Value2 = -371.8571999999999
$ws.Range("N74").Value2 = -19354.334
$ws.Range("H77").Value2 = 3290.9167
$ws.Range("I77").Value2 = 1245.8572
$ws.Range("J77").Value2 = 17606.334
$ws.Range("K77").Value2 = 6229.286
$ws.Range("L77").Value2 = 88031.67
$ws.Range("M77").Value2 = -1861.286
$ws.Range("N77").Value2 = -96767.67
$ws.Range("H93").Value2 = 16384
$ws.Range("J93").Value2 = 16384
$ws.Range("L93").Value2 = 16384
$ws.Range("N93").Value2 = -21376
$ws.Range("H95").Value2 = 30272.428
$ws.Range("J95").Value2 = 30272.428
$ws.Range("L95").Value2 = 30272.428
$ws.Range("N95").Value2 = -35764.428
$ws.Range("H101").Value2 = 138877.33
$ws.Range("J101").Value2 = 138877.33
$ws.Range("L101").Value2 = 138877.33
$ws.Range("N101").Value2 = -145367.33
$ws.Range("H103").Value2 = 81500
$ws.Range("J103").Value2 = 81500
$ws.Range("L103").Value2 = 81500
$ws.Range("N103").Value2 = -83844
$ws.Range("H105").Value2 = 0
$ws.Range("J105").Value2 = 0
$ws.Range("L105").Value2 = 0
$ws.Range("N105").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value2 = 150406.25
$ws.Range("J63").Value2 = 150406.25
$ws.Range("L63").Value2 = 150406.25
$ws.Range("N63").Value2 = -151778.25
$ws.Range("H66").Value2 = 150406.25
$ws.Range("J66").Value2 = 150406.25
$ws.Range("L66").Value2 = 451218.75
$ws.Range("N66").Value2 = -458082.75
$ws.Range("H100").Value2 = 142924.67
$ws.Range("J100").Value2 = 142924.67
$ws.Range("L100").Value2 = 142924.67
$ws.Range("N100").Value2 = -145088.67

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 8439.637000000001
$ws.Range("I31").Value2 = 3503.7896
$ws.Range("J31").Value2 = 39700
$ws.Range("K31").Value2 = 3503.7896
$ws.Range("L31").Value2 = 39700
$ws.Range("M31").Value2 = -3208.7896
$ws.Range("N31").Value2 = -40290
$ws.Range("H34").Value2 = 8439.637000000001
$ws.Range("I34").Value2 = 3503.7896
$ws.Range("J34").Value2 = 39700
$ws.Range("K34").Value2 = 3503.7896
$ws.Range("L34").Value2 = 39700
$ws.Range("M34").Value2 = -3301.7896
$ws.Range("N34").Value2 = -40104
$ws.Range("H39").Value2 = 16082.4
$ws.Range("I39").Value2 = 14228.125
$ws.Range("K39").Value2 = 14228.125
$ws.Range("M39").Value2 = -13837.125
$ws.Range("H49").Value2 = 16082.4
$ws.Range("I49").Value2 = 14228.125
$ws.Range("K49").Value2 = 14228.125
$ws.Range("M49").Value2 = -14046.125
$ws.Range("H74").Value2 = 131813.27
$ws.Range("I74").Value2 = 45300
$ws.Range("K74").Value2 = 45300
$ws.Range("M74").Value2 = -44426
$ws.Range("H77").Value2 = 131813.27
$ws.Range("I77").Value2 = 45300
$ws.Range("K77").Value2 = 135900
$ws.Range("M77").Value2 = -131532
$ws.Range("H105").Value2 = 2357.1667
$ws.Range("I105").Value2 = 2937
$ws.Range("J105").Value2 = 1197.5
$ws.Range("K105").Value2 = 2937
$ws.Range("L105").Value2 = 1197.5
$ws.Range("M105").Value2 = -1190
$ws.Range("N105").Value2 = -4691.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value2 = 335
$ws.Range("I8").Value2 = 335
$ws.Range("K8").Value2 = 1005
$ws.Range("M8").Value2 = -866
$ws.Range("H105").Value2 = 11494
$ws.Range("J105").Value2 = 11494
$ws.Range("L105").Value2 = 34482
$ws.Range("N105").Value2 = -39724

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value2 = 146235.58
$ws.Range("J95").Value2 = 146235.58
$ws.Range("L95").Value2 = 146235.58
$ws.Range("N95").Value2 = -151727.58
$ws.Range("H126").Value2 = 5715.4287
$ws.Range("I126").Value2 = 3252.25
$ws.Range("J126").Value2 = 8999.666999999999
$ws.Range("K126").Value2 = 9756.75
$ws.Range("L126").Value2 = 26999.001
$ws.Range("M126").Value2 = -7286.75
$ws.Range("N126").Value2 = -31939.001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value2 = 5031.0884
$ws.Range("I136").Value2 = 4275.5
$ws.Range("K136").Value2 = 12826.5
$ws.Range("M136").Value2 = -10276.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value2 = 115310
$ws.Range("J95").Value2 = 115310
$ws.Range("L95").Value2 = 115310
$ws.Range("N95").Value2 = -120802
